$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "305.18"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.23%"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.29%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.049"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.22%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08011"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.49%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.871"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-2.51%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.786"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9234"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.54%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1288"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-7.54%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1903"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.35%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09132"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.11%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03421"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-4.75%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09857"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.45%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001413"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.40%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.006225"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "5.46%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.845"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "8.25%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.128"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.23%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.353"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "16.34%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3418"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.45%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "3.39%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.818"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.55%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2308"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-7.97%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04427"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.23%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.98%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004886"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2.14%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-24.18%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01935"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-1.31%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05168"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "5.61%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007554"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.17%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01016"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "9.84%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1350"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.60%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002172"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "3.42%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009631"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-15.04%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006198"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-2.70%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.04%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "65.14"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "2.48%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "39.39%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.04%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.04%"
